$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: replace RAFAEL AUGUSTO BARRAZA RUIZ's record with
# KEIDID MERCEDES TORRES ACOSTA's record
$ws.Range("C16").Value = "1101457308"
$ws.Range("D16").Value = "KEIDID MERCEDES TORRES ACOSTA"
$ws.Range("E16").Value = "1803"
$ws.Range("F16").Value = 31249
$ws.Range("G16").Value = 781242

# Row 19: replace KEIDID MERCEDES TORRES ACOSTA's record with
# RAFAEL AUGUSTO BARRAZA RUIZ's record (with an updated Valor Mora)
$ws.Range("C19").Value = "73582352"
$ws.Range("D19").Value = "RAFAEL AUGUSTO BARRAZA RUIZ"
$ws.Range("E19").Value = "2208"
$ws.Range("F19").Value = 18666
$ws.Range("G19").Value = 1792700

$wb.Save()
